$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1302
$ws1.Range("F6").Value = 18188
$ws1.Range("F7").Value = 367
$ws1.Range("F9").Value = 1068
$ws1.Range("F10").Value = 6858
$ws1.Range("F12").Value = 160
$ws1.Range("F14").Value = 112
$ws1.Range("F17").Value = 158
$ws1.Range("F19").Value = 235
$ws1.Range("F20").Value = 57
$ws1.Range("F25").Value = 277
$ws1.Range("F26").Value = 990
$ws1.Range("F30").Value = 36
$ws1.Range("F32").Value = 73
$ws1.Range("F33").Value = 12085
$ws1.Range("F34").Value = 1282
$ws1.Range("F36").Value = 209
$ws1.Range("F37").Value = 284
$ws1.Range("F38").Value = 3921
$ws1.Range("F39").Value = 301

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1302
$ws4.Range("F6").Value = 18188
$ws4.Range("F7").Value = 367
$ws4.Range("F9").Value = 1068
$ws4.Range("F10").Value = 6858
$ws4.Range("F12").Value = 160
$ws4.Range("F14").Value = 112
$ws4.Range("F17").Value = 158
$ws4.Range("F19").Value = 235
$ws4.Range("F20").Value = 57
$ws4.Range("F25").Value = 277
$ws4.Range("F26").Value = 990
$ws4.Range("F32").Value = 36
$ws4.Range("F34").Value = 73
$ws4.Range("F35").Value = 12085
$ws4.Range("F36").Value = 1282
$ws4.Range("F38").Value = 209
$ws4.Range("F39").Value = 284
$ws4.Range("F40").Value = 3921
$ws4.Range("F41").Value = 301
